# Weekly price-sheet refresh: a new daily price record was added for the
# market "Terminal La Palmera de La Serena" / "Poroto granado".
#
# The new record belongs right after the existing row 10 (row 11), so every
# row that was previously 11..49 shifts down one position (becoming 12..50).
# Insert a whole row at 11 (this also pushes the sheet's used range from
# R49 to R50) and then fill that fresh row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Terminal La Palmera de La Serena"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44558
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 100112030
$ws.Range("G11").Value = "Poroto granado"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 400
$ws.Range("K11").Value = 30000
$ws.Range("L11").Value = 31000
$ws.Range("M11").Value = 30500
$ws.Range("N11").Value = "$/malla 25 kilos"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 1220
$ws.Range("Q11").Value = 25
$ws.Range("R11").Value = "Hortaliza"
